# Delete the first column ("id") from Hoja1, shifting remaining columns left.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").EntireColumn.Delete()
